# Apply symbol-list update (Tue Jan  3 08:16:34 UTC 2023 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''246.21'
$ws.Range('E2').Value = '''-0.09%'
$ws.Range('G2').Value = '''8'
$ws.Range('D3').Value = '''29.85'
$ws.Range('E3').Value = '''-0.21%'
$ws.Range('G3').Value = '''8'
$ws.Range('D4').Value = '''5.156'
$ws.Range('E4').Value = '''-0.13%'
$ws.Range('G4').Value = '''8'
$ws.Range('D5').Value = '''0.05789'
$ws.Range('E5').Value = '''1.06%'
$ws.Range('G5').Value = '''8'
$ws.Range('D6').Value = '''6.650'
$ws.Range('E6').Value = '''1.05%'
$ws.Range('G6').Value = '''8'
$ws.Range('D7').Value = '''3.229'
$ws.Range('E7').Value = '''7.08%'
$ws.Range('G7').Value = '''8'
$ws.Range('D8').Value = '''0.8498'
$ws.Range('E8').Value = '''-0.77%'
$ws.Range('G8').Value = '''8'
$ws.Range('D9').Value = '''0.8593'
$ws.Range('E9').Value = '''-1.13%'
$ws.Range('G9').Value = '''8'
$ws.Range('D10').Value = '''0.1383'
$ws.Range('E10').Value = '''2.19%'
$ws.Range('G10').Value = '''8'
$ws.Range('D11').Value = '''0.07104'
$ws.Range('E11').Value = '''1.85%'
$ws.Range('G11').Value = '''8'
$ws.Range('D12').Value = '''0.03252'
$ws.Range('E12').Value = '''11.35%'
$ws.Range('G12').Value = '''8'
$ws.Range('D13').Value = '''0.09371'
$ws.Range('E13').Value = '''0.05%'
$ws.Range('G13').Value = '''8'
$ws.Range('D14').Value = '''0.001538'
$ws.Range('E14').Value = '''0.85%'
$ws.Range('G14').Value = '''8'
$ws.Range('D15').Value = '''0.0006031'
$ws.Range('E15').Value = '''0.28%'
$ws.Range('G15').Value = '''8'
$ws.Range('D16').Value = '''0.006062'
$ws.Range('E16').Value = '''1.05%'
$ws.Range('G16').Value = '''8'
$ws.Range('D17').Value = '''3.505'
$ws.Range('E17').Value = '''-0.10%'
$ws.Range('G17').Value = '''8'
$ws.Range('D18').Value = '''2.211'
$ws.Range('E18').Value = '''1.31%'
$ws.Range('G18').Value = '''8'
$ws.Range('D19').Value = '''0.3165'
$ws.Range('E19').Value = '''0.67%'
$ws.Range('G19').Value = '''8'
$ws.Range('D20').Value = '''0.03367'
$ws.Range('E20').Value = '''1.64%'
$ws.Range('G20').Value = '''8'
$ws.Range('D21').Value = '''0.1297'
$ws.Range('E21').Value = '''-0.57%'
$ws.Range('G21').Value = '''8'
$ws.Range('D22').Value = '''3.491'
$ws.Range('E22').Value = '''-2.90%'
$ws.Range('G22').Value = '''8'
$ws.Range('D23').Value = '''0.04136'
$ws.Range('E23').Value = '''-0.31%'
$ws.Range('G23').Value = '''8'
$ws.Range('D24').Value = '''0.1381'
$ws.Range('E24').Value = '''0.38%'
$ws.Range('G24').Value = '''8'
$ws.Range('D25').Value = '''0.001228'
$ws.Range('E25').Value = '''1.36%'
$ws.Range('G25').Value = '''8'
$ws.Range('D26').Value = '''0.004141'
$ws.Range('E26').Value = '''-7.90%'
$ws.Range('G26').Value = '''8'
$ws.Range('E27').Value = '''1.98%'
$ws.Range('G27').Value = '''8'
$ws.Range('D28').Value = '''0.0001449'
$ws.Range('E28').Value = '''100.05%'
$ws.Range('G28').Value = '''8'
$ws.Range('G29').Value = '''8'
$ws.Range('G30').Value = '''8'
$ws.Range('G31').Value = '''8'
$ws.Range('G32').Value = '''8'
$ws.Range('G33').Value = '''8'
$ws.Range('G34').Value = '''8'
$ws.Range('G35').Value = '''8'
$ws.Range('G36').Value = '''8'
$ws.Range('G37').Value = '''8'
$ws.Range('G38').Value = '''8'
$ws.Range('G39').Value = '''8'
$ws.Range('D40').Value = '''0.03765'
$ws.Range('E40').Value = '''-0.58%'
$ws.Range('G40').Value = '''8'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '''0.1070'
$ws.Range('E41').Value = '''0.28%'
$ws.Range('G41').Value = '''8'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '''0.002470'
$ws.Range('E42').Value = '''-3.62%'
$ws.Range('G42').Value = '''8'
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = '''0.003528'
$ws.Range('E43').Value = '''-38.26%'
$ws.Range('G43').Value = '''8'
$ws.Range('E44').Value = '''-4.99%'
$ws.Range('G44').Value = '''8'
$ws.Range('D45').Value = '''0.00005323'
$ws.Range('E45').Value = '''4.55%'
$ws.Range('G45').Value = '''8'
$ws.Range('E46').Value = '''0.28%'
$ws.Range('G46').Value = '''8'
$ws.Range('D47').Value = '''0.07101'
$ws.Range('E47').Value = '''-11.00%'
$ws.Range('G47').Value = '''8'
$ws.Range('D48').Value = '''0.002184'
$ws.Range('E48').Value = '''-19.88%'
$ws.Range('G48').Value = '''8'
$ws.Range('D49').Value = '''0.00002100'
$ws.Range('E49').Value = '''0.28%'
$ws.Range('G49').Value = '''8'
$ws.Range('D50').Value = '''0.0002000'
$ws.Range('E50').Value = '''0.28%'
$ws.Range('G50').Value = '''8'
$ws.Range('G51').Value = '''8'
